$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the current price input ("Precio actual")
$ws.Range("G3").Value = 223.88

# Update the "Min P/E 10 years" input used to estimate future market price
$ws.Range("F12").Value = 60

# Add a new (currently blank) formatted cell below the warning box -
# a placeholder used to get the estimations for the stocks.
$ws.Range("G23").Value = ""
$ws.Range("G23").Font.Name = "Arial"
$ws.Range("G23").Font.Size = 10
$ws.Range("G23").Font.Underline = [int]2
$ws.Range("G23").Font.Color = 0
$ws.Rows.Item(23).RowHeight = 12.8

# Move the active selection to F13
$ws.Range("F13").Select() | Out-Null
